# Update Name of Algo
# Apply updated numeric values to the result_data_KNN sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 16.359
$ws.Range("C12").Value = -10.94
$ws.Range("E14").Value = 17.009
$ws.Range("E26").Value = 16.407
$ws.Range("C27").Value = -13.058
$ws.Range("E31").Value = 16.331
$ws.Range("C32").Value = -13.18
$ws.Range("E35").Value = 16.492
$ws.Range("C36").Value = -12.732
$ws.Range("E37").Value = 16.815
$ws.Range("C38").Value = -12.607
$ws.Range("E45").Value = 16.662
$ws.Range("C46").Value = -13.742
$ws.Range("E52").Value = 17.138
$ws.Range("C54").Value = -13.137
$ws.Range("C55").Value = -13.46
$ws.Range("C56").Value = -13.364
$ws.Range("E57").Value = 16.45
$ws.Range("C67").Value = -12.037
$ws.Range("C69").Value = -11.062
$ws.Range("C72").Value = -11.932
$ws.Range("E81").Value = 16.517
$ws.Range("C83").Value = -13.213
$ws.Range("E83").Value = 16.709
$ws.Range("C86").Value = -13.967
$ws.Range("C91").Value = -11.444
$ws.Range("C93").Value = -11.979
$ws.Range("C99").Value = -12.635
$ws.Range("E100").Value = 16.555
$ws.Range("E102").Value = 16.519
